$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 2979.2964
$ws.Cells.Item(51, 9).Value = 1171.4166
$ws.Cells.Item(51, 10).Value = 4425.6
$ws.Cells.Item(51, 11).Value = 1171.4166
$ws.Cells.Item(51, 12).Value = 4425.6
$ws.Cells.Item(51, 13).Value = -687.4166
$ws.Cells.Item(51, 14).Value = -5393.6
$ws.Cells.Item(62, 8).Value = 2080.5715
$ws.Cells.Item(62, 9).Value = 2224.3333
$ws.Cells.Item(62, 10).Value = 1888.8889
$ws.Cells.Item(62, 11).Value = 2224.3333
$ws.Cells.Item(62, 12).Value = 1888.8889
$ws.Cells.Item(62, 13).Value = -1600.3333
$ws.Cells.Item(62, 14).Value = -3136.8889
$ws.Cells.Item(65, 8).Value = 2080.5715
$ws.Cells.Item(65, 9).Value = 2224.3333
$ws.Cells.Item(65, 10).Value = 1888.8889
$ws.Cells.Item(65, 11).Value = 11121.6665
$ws.Cells.Item(65, 12).Value = 9444.4445
$ws.Cells.Item(65, 13).Value = -8001.666499999999
$ws.Cells.Item(65, 14).Value = -15684.4445
$ws.Cells.Item(76, 8).Value = 3370.2778
$ws.Cells.Item(76, 9).Value = 3364.2144
$ws.Cells.Item(76, 10).Value = 3391.5
$ws.Cells.Item(76, 11).Value = 3364.2144
$ws.Cells.Item(76, 12).Value = 3391.5
$ws.Cells.Item(76, 13).Value = -3049.2144
$ws.Cells.Item(76, 14).Value = -4021.5
$ws.Cells.Item(79, 8).Value = 3370.2778
$ws.Cells.Item(79, 9).Value = 3364.2144
$ws.Cells.Item(79, 10).Value = 3391.5
$ws.Cells.Item(79, 11).Value = 3364.2144
$ws.Cells.Item(79, 12).Value = 3391.5
$ws.Cells.Item(79, 13).Value = -2272.2144
$ws.Cells.Item(79, 14).Value = -5575.5
$ws.Cells.Item(106, 8).Value = 6252.8823
$ws.Cells.Item(106, 9).Value = 1922.9231
$ws.Cells.Item(106, 10).Value = 8933.333000000001
$ws.Cells.Item(106, 11).Value = 1922.9231
$ws.Cells.Item(106, 12).Value = 8933.333000000001
$ws.Cells.Item(106, 13).Value = -1291.9231
$ws.Cells.Item(106, 14).Value = -10195.333
$ws.Cells.Item(113, 8).Value = 4142.846
$ws.Cells.Item(113, 9).Value = 4044.625
$ws.Cells.Item(113, 10).Value = 4300
$ws.Cells.Item(113, 11).Value = 4044.625
$ws.Cells.Item(113, 12).Value = 4300
$ws.Cells.Item(113, 13).Value = -790.625
$ws.Cells.Item(113, 14).Value = -10808
$ws.Cells.Item(138, 8).Value = 4446619
$ws.Cells.Item(138, 9).Value = 1890.875
$ws.Cells.Item(138, 10).Value = 6898883
$ws.Cells.Item(138, 11).Value = 5672.625
$ws.Cells.Item(138, 12).Value = 20696649
$ws.Cells.Item(138, 13).Value = -532.625
$ws.Cells.Item(138, 14).Value = -20706929
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(92, 8).Value = 31987.5
$ws.Cells.Item(92, 10).Value = 31987.5
$ws.Cells.Item(92, 12).Value = 31987.5
$ws.Cells.Item(92, 14).Value = -36979.5
$ws.Cells.Item(122, 8).Value = 1813.6
$ws.Cells.Item(122, 9).Value = 1289
$ws.Cells.Item(122, 10).Value = 2862.8
$ws.Cells.Item(122, 11).Value = 3867
$ws.Cells.Item(122, 12).Value = 8588.400000000001
$ws.Cells.Item(122, 13).Value = -1417
$ws.Cells.Item(122, 14).Value = -13488.4
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(76, 8).Value = 0
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).Value = ""
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).Value = ""
$ws.Cells.Item(86, 8).Value = 13150.143
$ws.Cells.Item(86, 9).Value = 15132.117
$ws.Cells.Item(86, 10).Value = 4726.75
$ws.Cells.Item(86, 11).Value = 15132.117
$ws.Cells.Item(86, 12).Value = 4726.75
$ws.Cells.Item(86, 13).Value = -14009.117
$ws.Cells.Item(86, 14).Value = -6972.75
$ws.Cells.Item(89, 8).Value = 13150.143
$ws.Cells.Item(89, 9).Value = 15132.117
$ws.Cells.Item(89, 10).Value = 4726.75
$ws.Cells.Item(89, 11).Value = 75660.58500000001
$ws.Cells.Item(89, 12).Value = 23633.75
$ws.Cells.Item(89, 13).Value = -70044.58500000001
$ws.Cells.Item(89, 14).Value = -34865.75
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 218287.53
$ws.Cells.Item(31, 9).Value = 47618.637
$ws.Cells.Item(31, 10).Value = 328720.34
$ws.Cells.Item(31, 11).Value = 47618.637
$ws.Cells.Item(31, 12).Value = 328720.34
$ws.Cells.Item(31, 13).Value = -47323.637
$ws.Cells.Item(31, 14).Value = -329310.34
$ws.Cells.Item(34, 8).Value = 218287.53
$ws.Cells.Item(34, 9).Value = 47618.637
$ws.Cells.Item(34, 10).Value = 328720.34
$ws.Cells.Item(34, 11).Value = 47618.637
$ws.Cells.Item(34, 12).Value = 328720.34
$ws.Cells.Item(34, 13).Value = -47416.637
$ws.Cells.Item(34, 14).Value = -329124.34
$ws.Cells.Item(41, 8).Value = 11050
$ws.Cells.Item(41, 9).Value = 7100
$ws.Cells.Item(41, 10).Value = 15000
$ws.Cells.Item(41, 11).Value = 7100
$ws.Cells.Item(41, 12).Value = 15000
$ws.Cells.Item(41, 13).Value = -6672
$ws.Cells.Item(41, 14).Value = -15856
$ws.Cells.Item(42, 8).Value = 0
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 11).Value = 0
$ws.Cells.Item(42, 13).Value = ""
$ws.Cells.Item(95, 8).Value = 30174.666
$ws.Cells.Item(95, 10).Value = 30174.666
$ws.Cells.Item(95, 12).Value = 30174.666
$ws.Cells.Item(95, 14).Value = -35666.666
$ws.Cells.Item(134, 8).Value = 48588.19
$ws.Cells.Item(134, 9).Value = 855.63635
$ws.Cells.Item(134, 10).Value = 101094
$ws.Cells.Item(134, 11).Value = 2566.90905
$ws.Cells.Item(134, 12).Value = 303282
$ws.Cells.Item(134, 13).Value = -31.90905000000021
$ws.Cells.Item(134, 14).Value = -308352
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 256.60715
$ws.Cells.Item(33, 9).Value = 155.86957
$ws.Cells.Item(33, 11).Value = 935.2174200000001
$ws.Cells.Item(33, 13).Value = -652.2174200000001
$ws.Cells.Item(76, 8).Value = 3300
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 11).Value = 0
$ws.Cells.Item(76, 13).Value = ""
$ws.Cells.Item(79, 8).Value = 3300
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 11).Value = 0
$ws.Cells.Item(79, 13).Value = ""
$ws.Cells.Item(131, 8).Value = 969
$ws.Cells.Item(131, 9).Value = 390
$ws.Cells.Item(131, 10).Value = 1068.8276
$ws.Cells.Item(131, 11).Value = 1170
$ws.Cells.Item(131, 12).Value = 3206.4828
$ws.Cells.Item(131, 13).Value = 3870
$ws.Cells.Item(131, 14).Value = -13286.4828
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2592
$ws.Cells.Item(40, 9).Value = 2024
$ws.Cells.Item(40, 10).Value = 6000
$ws.Cells.Item(40, 11).Value = 2024
$ws.Cells.Item(40, 12).Value = 6000
$ws.Cells.Item(40, 13).Value = -1888
$ws.Cells.Item(40, 14).Value = -6272
$ws.Cells.Item(82, 8).Value = 1529.2222
$ws.Cells.Item(82, 9).Value = 1247.1666
$ws.Cells.Item(82, 11).Value = 1247.1666
$ws.Cells.Item(82, 13).Value = -886.1666
$ws.Cells.Item(85, 8).Value = 1529.2222
$ws.Cells.Item(85, 9).Value = 1247.1666
$ws.Cells.Item(85, 11).Value = 1247.1666
$ws.Cells.Item(85, 13).Value = 0.8333999999999833
$ws.Cells.Item(132, 8).Value = 35153
$ws.Cells.Item(132, 9).Value = 2664.5417
$ws.Cells.Item(132, 10).Value = 146542
$ws.Cells.Item(132, 11).Value = 7993.625100000001
$ws.Cells.Item(132, 12).Value = 439626
$ws.Cells.Item(132, 13).Value = -5463.625100000001
$ws.Cells.Item(132, 14).Value = -444686
$ws.Cells.Item(136, 8).Value = 66354.03
$ws.Cells.Item(136, 9).Value = 42321.117
$ws.Cells.Item(136, 10).Value = 170496.67
$ws.Cells.Item(136, 11).Value = 126963.351
$ws.Cells.Item(136, 12).Value = 511490.01
$ws.Cells.Item(136, 13).Value = -124413.351
$ws.Cells.Item(136, 14).Value = -516590.01
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(49, 8).Value = 15733.75
$ws.Cells.Item(49, 10).Value = 15733.75
$ws.Cells.Item(49, 12).Value = 15733.75
$ws.Cells.Item(49, 14).Value = -16193.75
$ws.Cells.Item(100, 8).Value = 144943
$ws.Cells.Item(100, 9).Value = 167534
$ws.Cells.Item(100, 10).Value = 127999.75
$ws.Cells.Item(100, 11).Value = 335068
$ws.Cells.Item(100, 12).Value = 255999.5
$ws.Cells.Item(100, 13).Value = -334527
$ws.Cells.Item(100, 14).Value = -257081.5
$ws.Cells.Item(132, 8).Value = 92360.37
$ws.Cells.Item(132, 9).Value = 63433.125
$ws.Cells.Item(132, 10).Value = 169499.67
$ws.Cells.Item(132, 11).Value = 190299.375
$ws.Cells.Item(132, 12).Value = 508499.01
$ws.Cells.Item(132, 13).Value = -187769.375
$ws.Cells.Item(132, 14).Value = -513559.01
